$d = $word.ActiveDocument

# --- 1. Merge the split runs in the "La classe Scheduler ..." paragraph ---
$pScheduler = $d.Paragraphs.Item(54)
$rScheduler = $pScheduler.Range
$rScheduler.Find.Execute("La class", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "La class", 2)

# --- 2. Merge the split runs in the "... collection de lois de probabilité ..." paragraph ---
$pStdRandom = $d.Paragraphs.Item(56)
$rStdRandom = $pStdRandom.Range
$rStdRandom.Find.Execute("collection de loi", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "collection de loi", 2)

# --- 3. Append new content after the last paragraph ("Ces deux fonctions ... aléatoires.") ---

# empty paragraph
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $pLast.Range
$rLast.Collapse(0)
$rLast.InsertParagraphAfter()

# heading paragraph: "Visualisation du graphe des recettes" (underlined)
$pEmpty = $d.Paragraphs.Item($d.Paragraphs.Count)
$rEmpty = $pEmpty.Range
$rEmpty.Collapse(0)
$rEmpty.InsertParagraphAfter()
$pHeading = $d.Paragraphs.Item($d.Paragraphs.Count)
$rHeading = $pHeading.Range
$rHeading.Font.Underline = 1
$rHeading.InsertAfter("Visualisation du graphe des recettes")

# paragraph: Gephi intro text (underline none)
$pHeading2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$rHeading2 = $pHeading2.Range
$rHeading2.Collapse(0)
$rHeading2.InsertParagraphAfter()
$pGephi = $d.Paragraphs.Item($d.Paragraphs.Count)
$rGephi = $pGephi.Range
$rGephi.Font.Underline = 0
$rGephi.InsertAfter("Puisque les recettes lie un acteur à un autre, il est intéressant de les visualiser dans un graphe. Nous utilisons pour cela l'outil de gestion de graphe Gephi. Puisque la génération du graphe est complètement indépendante de la simulation et prend un peu de temps, nous avons choisi de la traiter dans un exécutable ajouté en tant que bibliothèque, ce qui présente l'avantage d'un faible couplage.")

# empty paragraph (underline none)
$pGephi2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$rGephi2 = $pGephi2.Range
$rGephi2.Collapse(0)
$rGephi2.InsertParagraphAfter()
$pBlank2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$rBlank2 = $pBlank2.Range
$rBlank2.Font.Underline = 0

# paragraph: export/graph details text (underline none)
$rBlank2b = $pBlank2.Range
$rBlank2b.Collapse(0)
$rBlank2b.InsertParagraphAfter()
$pDetails = $d.Paragraphs.Item($d.Paragraphs.Count)
$rDetails = $pDetails.Range
$rDetails.Font.Underline = 0
$rDetails.InsertAfter("Sans entrer dans les détails techniques de ce qui n'est pas le cœur de notre simulation, on peut dire que la communication s'effectue par fichiers : quand l'ensemble des données change, un fichier d'export est généré. Suivant les arguments qu'on lui passe, l'exécutable lit ce fichier et produit un graphe. Deux supports de productions sont possibles : un fichier svg et un applet Java. La deuxième possibilité offre la possibilité d'explorer le graphe. Enfin, le fichier d'export généré peut aussi être importé dans Gephi lui-même pour des calculs plus poussés sur ce graphe (rayon, diamètre et autre).")

# --- 4. styles.xml: remove default run size (w:sz) from rPrDefault and set Normal style color ---
$styles = $d.Styles
$normalStyle = $styles.Item("Normal")
$normalStyle.Font.Color = 655360

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
